$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.07173050923138437
$ws.Range("J2").Value = 0.07173050923138438
$ws.Range("M2").Value = 6.369648000000001
$ws.Range("N2").Value = 19.108944
$ws.Range("O2").Value = 0.1127004548956141
$ws.Range("P2").Value = 0.1127004548956141
$ws.Range("Q2").Value = 0.19578174736
$ws.Range("R2").Value = 1.76203572624
$ws.Range("S2").Value = 0.008084061020271066
$ws.Range("T2").Value = 0.008084061020271068

$ws.Range("I3").Value = 0.07173050923138437
$ws.Range("J3").Value = 0.07173050923138438
$ws.Range("O3").Value = 0.2159356303411415
$ws.Range("P3").Value = 0.2159356303411415
$ws.Range("S3").Value = 0.01548917272557005
$ws.Range("T3").Value = 0.01548917272557006

$ws.Range("I4").Value = 0.07173050923138437
$ws.Range("J4").Value = 0.07173050923138438
$ws.Range("M4").Value = 5.679255
$ws.Range("N4").Value = 17.037765
$ws.Range("O4").Value = 0.1004850851990865
$ws.Range("P4").Value = 0.1004850851990865
$ws.Range("Q4").Value = 0.17456136785
$ws.Range("R4").Value = 1.57105231065
$ws.Range("S4").Value = 0.00720784633148952
$ws.Range("T4").Value = 0.007207846331489521

$ws.Range("I5").Value = 0.07173050923138437
$ws.Range("J5").Value = 0.07173050923138438
$ws.Range("M5").Value = 32.265151
$ws.Range("N5").Value = 96.79545300000001
$ws.Range("O5").Value = 0.5708788295641578
$ws.Range("P5").Value = 0.5708788295641579
$ws.Range("Q5").Value = 0.9917231912366667
$ws.Range("R5").Value = 8.92550872113
$ws.Range("S5").Value = 0.04094942915405372
$ws.Range("T5").Value = 0.04094942915405374

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3977653333333333
$ws.Range("H6").Value = 1.193296
$ws.Range("I6").Value = 0.9282694907686155
$ws.Range("J6").Value = 0.9282694907686156
$ws.Range("M6").Value = 6.369648000000001
$ws.Range("N6").Value = 19.108944
$ws.Range("O6").Value = 0.1127004548956141
$ws.Range("P6").Value = 0.1127004548956141
$ws.Range("Q6").Value = 2.533625159936
$ws.Range("R6").Value = 22.802626439424
$ws.Range("S6").Value = 0.1046163938753431
$ws.Range("T6").Value = 0.1046163938753431

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3977653333333333
$ws.Range("H7").Value = 1.193296
$ws.Range("I7").Value = 0.9282694907686155
$ws.Range("J7").Value = 0.9282694907686156
$ws.Range("O7").Value = 0.2159356303411415
$ws.Range("P7").Value = 0.2159356303411415
$ws.Range("Q7").Value = 4.854460849032888
$ws.Range("R7").Value = 43.69014764129599
$ws.Range("S7").Value = 0.2004464576155714
$ws.Range("T7").Value = 0.2004464576155715

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3977653333333333
$ws.Range("H8").Value = 1.193296
$ws.Range("I8").Value = 0.9282694907686155
$ws.Range("J8").Value = 0.9282694907686156
$ws.Range("M8").Value = 5.679255
$ws.Range("N8").Value = 17.037765
$ws.Range("O8").Value = 0.1004850851990865
$ws.Range("P8").Value = 0.1004850851990865
$ws.Range("Q8").Value = 2.25901075816
$ws.Range("R8").Value = 20.33109682344
$ws.Range("S8").Value = 0.09327723886759698
$ws.Range("T8").Value = 0.09327723886759699

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3977653333333333
$ws.Range("H9").Value = 1.193296
$ws.Range("I9").Value = 0.9282694907686155
$ws.Range("J9").Value = 0.9282694907686156
$ws.Range("M9").Value = 32.265151
$ws.Range("N9").Value = 96.79545300000001
$ws.Range("O9").Value = 0.5708788295641578
$ws.Range("P9").Value = 0.5708788295641579
$ws.Range("Q9").Value = 12.83395854256533
$ws.Range("R9").Value = 115.505626883088
$ws.Range("S9").Value = 0.5299294004101041
$ws.Range("T9").Value = 0.5299294004101042
